$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.659013
$ws.Cells.Item(2, 8).Value = 7.977039
$ws.Cells.Item(2, 9).Value = 0.4672138103908383
$ws.Cells.Item(2, 10).Value = 0.4672138103908383
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.664704
$ws.Cells.Item(2, 14).Value = 3.329408
$ws.Cells.Item(2, 15).Value = 0.03573178976417293
$ws.Cells.Item(2, 16).Value = 0.02637436895330882
$ws.Cells.Item(2, 17).Value = 4.426469577152
$ws.Cells.Item(2, 18).Value = 26.558817462912
$ws.Cells.Item(2, 19).Value = 0.01669438564780359
$ws.Cells.Item(2, 20).Value = 0.01232246941532924

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.659013
$ws.Cells.Item(3, 8).Value = 7.977039
$ws.Cells.Item(3, 9).Value = 0.4672138103908383
$ws.Cells.Item(3, 10).Value = 0.4672138103908383
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 26.69016266666667
$ws.Cells.Item(3, 14).Value = 80.070488
$ws.Cells.Item(3, 15).Value = 0.572887000437862
$ws.Cells.Item(3, 16).Value = 0.6342895171704658
$ws.Cells.Item(3, 17).Value = 70.96948950278133
$ws.Cells.Item(3, 18).Value = 638.725405525032
$ws.Cells.Item(3, 19).Value = 0.2676607183979514
$ws.Cells.Item(3, 20).Value = 0.2963488222081783

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.659013
$ws.Cells.Item(4, 8).Value = 7.977039
$ws.Cells.Item(4, 9).Value = 0.4672138103908383
$ws.Cells.Item(4, 10).Value = 0.4672138103908383
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.067451333333333
$ws.Cells.Item(4, 14).Value = 3.202354
$ws.Cells.Item(4, 15).Value = 0.0229121493227341
$ws.Cells.Item(4, 16).Value = 0.02536789300533437
$ws.Cells.Item(4, 17).Value = 2.838366972200666
$ws.Cells.Item(4, 18).Value = 25.54530274980599
$ws.Cells.Item(4, 19).Value = 0.01070487258931846
$ws.Cells.Item(4, 20).Value = 0.01185222995260937

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.659013
$ws.Cells.Item(5, 8).Value = 7.977039
$ws.Cells.Item(5, 9).Value = 0.4672138103908383
$ws.Cells.Item(5, 10).Value = 0.4672138103908383
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.531414666666667
$ws.Cells.Item(5, 14).Value = 7.594244
$ws.Cells.Item(5, 15).Value = 0.05433517110265684
$ws.Cells.Item(5, 16).Value = 0.06015886102798208
$ws.Cells.Item(5, 17).Value = 6.731064507057333
$ws.Cells.Item(5, 18).Value = 60.579580563516
$ws.Cells.Item(5, 19).Value = 0.02538614232911047
$ws.Cells.Item(5, 20).Value = 0.02810705068965641

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.659013
$ws.Cells.Item(6, 8).Value = 7.977039
$ws.Cells.Item(6, 9).Value = 0.4672138103908383
$ws.Cells.Item(6, 10).Value = 0.4672138103908383
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.769714
$ws.Cells.Item(6, 14).Value = 8.309142000000001
$ws.Cells.Item(6, 15).Value = 0.05945011146419213
$ws.Cells.Item(6, 16).Value = 0.06582202505473476
$ws.Cells.Item(6, 17).Value = 7.364705532282001
$ws.Cells.Item(6, 18).Value = 66.282349790538
$ws.Cells.Item(6, 19).Value = 0.02777591310534526
$ws.Cells.Item(6, 20).Value = 0.03075295913346385

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.659013
$ws.Cells.Item(7, 8).Value = 7.977039
$ws.Cells.Item(7, 9).Value = 0.4672138103908383
$ws.Cells.Item(7, 10).Value = 0.4672138103908383
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 11.8654315
$ws.Cells.Item(7, 14).Value = 23.730863
$ws.Cells.Item(7, 15).Value = 0.254683777908382
$ws.Cells.Item(7, 16).Value = 0.187987334788174
$ws.Cells.Item(7, 17).Value = 31.5503366091095
$ws.Cells.Item(7, 18).Value = 189.302019654657
$ws.Cells.Item(7, 19).Value = 0.1189917783213091
$ws.Cells.Item(7, 20).Value = 0.08783027899160098

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.032199333333333
$ws.Cells.Item(8, 8).Value = 9.096598
$ws.Cells.Item(8, 9).Value = 0.5327861896091618
$ws.Cells.Item(8, 10).Value = 0.5327861896091618
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.664704
$ws.Cells.Item(8, 14).Value = 3.329408
$ws.Cells.Item(8, 15).Value = 0.03573178976417293
$ws.Cells.Item(8, 16).Value = 0.02637436895330882
$ws.Cells.Item(8, 17).Value = 5.047714358997333
$ws.Cells.Item(8, 18).Value = 30.286286153984
$ws.Cells.Item(8, 19).Value = 0.01903740411636935
$ws.Cells.Item(8, 20).Value = 0.01405189953797958

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.032199333333333
$ws.Cells.Item(9, 8).Value = 9.096598
$ws.Cells.Item(9, 9).Value = 0.5327861896091618
$ws.Cells.Item(9, 10).Value = 0.5327861896091618
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 26.69016266666667
$ws.Cells.Item(9, 14).Value = 80.070488
$ws.Cells.Item(9, 15).Value = 0.572887000437862
$ws.Cells.Item(9, 16).Value = 0.6342895171704658
$ws.Cells.Item(9, 17).Value = 80.92989344442488
$ws.Cells.Item(9, 18).Value = 728.3690409998239
$ws.Cells.Item(9, 19).Value = 0.3052262820399107
$ws.Cells.Item(9, 20).Value = 0.3379406949622875

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.032199333333333
$ws.Cells.Item(10, 8).Value = 9.096598
$ws.Cells.Item(10, 9).Value = 0.5327861896091618
$ws.Cells.Item(10, 10).Value = 0.5327861896091618
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.067451333333333
$ws.Cells.Item(10, 14).Value = 3.202354
$ws.Cells.Item(10, 15).Value = 0.0229121493227341
$ws.Cells.Item(10, 16).Value = 0.02536789300533437
$ws.Cells.Item(10, 17).Value = 3.236725221299111
$ws.Cells.Item(10, 18).Value = 29.130526991692
$ws.Cells.Item(10, 19).Value = 0.01220727673341564
$ws.Cells.Item(10, 20).Value = 0.01351566305272501

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 3.032199333333333
$ws.Cells.Item(11, 8).Value = 9.096598
$ws.Cells.Item(11, 9).Value = 0.5327861896091618
$ws.Cells.Item(11, 10).Value = 0.5327861896091618
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.531414666666667
$ws.Cells.Item(11, 14).Value = 7.594244
$ws.Cells.Item(11, 15).Value = 0.05433517110265684
$ws.Cells.Item(11, 16).Value = 0.06015886102798208
$ws.Cells.Item(11, 17).Value = 7.675753864656889
$ws.Cells.Item(11, 18).Value = 69.081784781912
$ws.Cells.Item(11, 19).Value = 0.02894902877354637
$ws.Cells.Item(11, 20).Value = 0.03205181033832567

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 3.032199333333333
$ws.Cells.Item(12, 8).Value = 9.096598
$ws.Cells.Item(12, 9).Value = 0.5327861896091618
$ws.Cells.Item(12, 10).Value = 0.5327861896091618
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.769714
$ws.Cells.Item(12, 14).Value = 8.309142000000001
$ws.Cells.Item(12, 15).Value = 0.05945011146419213
$ws.Cells.Item(12, 16).Value = 0.06582202505473476
$ws.Cells.Item(12, 17).Value = 8.398324944324001
$ws.Cells.Item(12, 18).Value = 75.58492449891601
$ws.Cells.Item(12, 19).Value = 0.03167419835884687
$ws.Cells.Item(12, 20).Value = 0.03506906592127092

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 3.032199333333333
$ws.Cells.Item(13, 8).Value = 9.096598
$ws.Cells.Item(13, 9).Value = 0.5327861896091618
$ws.Cells.Item(13, 10).Value = 0.5327861896091618
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 11.8654315
$ws.Cells.Item(13, 14).Value = 23.730863
$ws.Cells.Item(13, 15).Value = 0.254683777908382
$ws.Cells.Item(13, 16).Value = 0.187987334788174
$ws.Cells.Item(13, 17).Value = 35.97835348401233
$ws.Cells.Item(13, 18).Value = 215.870120904074
$ws.Cells.Item(13, 19).Value = 0.1356919995870728
$ws.Cells.Item(13, 20).Value = 0.1001570557965731

